# Refresh the crypto price/volume snapshot (Price = column D, Volume(1h) = column E)
# to match the latest scrape. Price cells that look like plain decimals (e.g. "9.04")
# are written with a leading apostrophe so Excel keeps them as text -- exactly like the
# existing inlineStr cells -- instead of silently re-typing them as numbers (which would
# also corrupt values such as "0.0480" by dropping the trailing zero).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '28.300.26'
$ws.Range("E2").Value = '  -0.74%  '
# Row 3: Ethereum
$ws.Range("D3").Value = '1.572.97'
$ws.Range("E3").Value = '  -0.36%  '
# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.01%  '
# Row 5: BNB
$ws.Range("D5").Value = '''211.97'
$ws.Range("E5").Value = '  -0.20%  '
# Row 6: XRP
$ws.Range("D6").Value = '''0.488'
$ws.Range("E6").Value = '  -0.69%  '
# Row 7: USDC
$ws.Range("E7").Value = '  +0.01%  '
# Row 8: OKB
$ws.Range("D8").Value = '''44.51'
$ws.Range("E8").Value = '  -5.43%  '
# Row 9: Solana
$ws.Range("D9").Value = '''23.74'
$ws.Range("E9").Value = '  -0.93%  '
# Row 10: Cardano
$ws.Range("D10").Value = '''0.246'
$ws.Range("E10").Value = '  -1.02%  '
# Row 11: Dogecoin
$ws.Range("E11").Value = '  -1.00%  '
# Row 12: TRON
$ws.Range("D12").Value = '''0.0895'
$ws.Range("E12").Value = '  +1.59%  '
# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = '1.797.02'
$ws.Range("E13").Value = '  -0.38%  '
# Row 14: WrappedEther
$ws.Range("D14").Value = '1.576.35'
$ws.Range("E14").Value = '  -0.13%  '
# Row 15: Polkadot
$ws.Range("D15").Value = '''3.69'
$ws.Range("E15").Value = '  -0.59%  '
# Row 16: Polygon
$ws.Range("D16").Value = '''0.516'
$ws.Range("E16").Value = '  -1.24%  '
# Row 17: WrappedBTC
$ws.Range("D17").Value = '28.321.05'
$ws.Range("E17").Value = '  -0.73%  '
# Row 18: Litecoin
$ws.Range("D18").Value = '''61.59'
$ws.Range("E18").Value = '  -1.38%  '
# Row 19: BitcoinCash
$ws.Range("D19").Value = '''229.91'
$ws.Range("E19").Value = '  +0.37%  '
# Row 20: Chainlink
$ws.Range("E20").Value = '  +0.01%  '
# Row 21: ShibaInu
$ws.Range("E21").Value = '  -1.39%  '
# Row 22: Dai
$ws.Range("E22").Value = '  -0.01%  '
# Row 23: Uniswap
$ws.Range("E23").Value = '  +0.46%  '
# Row 24: Avalanche
$ws.Range("D24").Value = '''9.04'
$ws.Range("E24").Value = '  -1.29%  '
# Row 25: Toncoin
$ws.Range("D25").Value = '''2.04'
$ws.Range("E25").Value = '  +0.35%  '
# Row 26: Monero
$ws.Range("D26").Value = '''151.45'
# Row 27: EthereumClassic
$ws.Range("D27").Value = '''14.94'
$ws.Range("E27").Value = '  -0.68%  '
# Row 28: Cosmos
$ws.Range("D28").Value = '''6.36'
$ws.Range("E28").Value = '  -1.71%  '
# Row 29: Stellar
$ws.Range("E29").Value = '  -1.98%  '
# Row 30: BinanceUSD
$ws.Range("E30").Value = '  -0.03%  '
# Row 31: Hedera
$ws.Range("D31").Value = '''0.0480'
$ws.Range("E31").Value = '  +2.95%  '
# Row 32: PancakeSwap
$ws.Range("E32").Value = '  -3.42%  '
# Row 33: Filecoin
$ws.Range("E33").Value = '  -0.62%  '
# Row 34: InternetComputer(DFINITY)
$ws.Range("D34").Value = '''3.08'
$ws.Range("E34").Value = '  -1.39%  '
# Row 35: Maker
$ws.Range("D35").Value = '1.380.91'
$ws.Range("E35").Value = '  -1.20%  '
# Row 36: TrustWalletToken
$ws.Range("E36").Value = '  +5.34%  '
# Row 37: LidoDAOToken
$ws.Range("E37").Value = '  -3.23%  '
# Row 38: HuobiToken
$ws.Range("D38").Value = '''2.37'
$ws.Range("E38").Value = '  +0.09%  '
# Row 39: MXToken
$ws.Range("E39").Value = '  +1.40%  '
# Row 40: VeChain
$ws.Range("E40").Value = '  -1.80%  '
# Row 41: ImmutableX
$ws.Range("E41").Value = '  -2.17%  '
# Row 42: PaxDollar
$ws.Range("E42").Value = '  -0.06%  '
# Row 43: RenderToken
$ws.Range("D43").Value = '''1.89'
$ws.Range("E43").Value = '  +1.86%  '
# Row 44: ARBITRUM
$ws.Range("E44").Value = '  -1.22%  '
# Row 45: Kaspa
$ws.Range("D45").Value = '''0.0463'
$ws.Range("E45").Value = '  -0.34%  '
# Row 46: FraxShare
$ws.Range("D46").Value = '''5.39'
$ws.Range("E46").Value = '  -3.78%  '
# Row 47: WEMIXToken
$ws.Range("D47").Value = '''0.923'
$ws.Range("E47").Value = '  -5.73%  '
# Row 48: Aave
$ws.Range("D48").Value = '''62.26'
$ws.Range("E48").Value = '  -0.66%  '
# Row 49: RocketPoolETH
$ws.Range("D49").Value = '1.709.75'
$ws.Range("E49").Value = '  -0.38%  '
# Row 50: mCoin
$ws.Range("E50").Value = '  +0.62%  '
# Row 51: Quant
$ws.Range("D51").Value = '''85.49'
$ws.Range("E51").Value = '  -0.42%  '
